$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix vendor name typo/formatting: "McMaster Carr" -> "McMaster-Carr" everywhere it appears.
$ws.Range("C11").Value = "McMaster-Carr"
$ws.Range("C12").Value = "McMaster-Carr"
$ws.Range("C13").Value = "McMaster-Carr"
$ws.Range("C14").Value = "McMaster-Carr"
$ws.Range("C15").Value = "McMaster-Carr"
$ws.Range("C16").Value = "McMaster-Carr"

# Add missing BOM item (row 17): item #16, a hex nut from McMaster-Carr.
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "Black-Oxide 18-8 Stainless Steel Hex Nut M2.5 x 0.45 mm Thread"
$ws.Range("C17").Value = "McMaster-Carr"
$ws.Range("D17").Value = "98676A320"
$ws.Range("E17").Value = 4

# F17 should take on the same (default) formatting as the other Product URL
# cells in this column (F2:F16), rather than the blank row's placeholder style.
$ws.Range("F16").Copy()
$ws.Range("F17").PasteSpecial(-4122)
$ws.Range("F17").Value = "https://www.mcmaster.com/98676A320/?SrchEntryWebPart_InpBox=98676A320"

$ws.Range("G17").Value = "NA"
$ws.Range("H17").Value = "NA"

# Update the saved cursor/selection position.
$ws.Range("F31").Select()
